$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($targetRange, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range("ZZ1").Formula = '="' + $escaped + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($targetRange).PasteSpecial(-4163)
}

Set-TextValue 'D2' '250.58'
Set-TextValue 'D3' '22.96'
Set-TextValue 'D4' '5.442'
Set-TextValue 'D5' '0.05677'
Set-TextValue 'D6' '3.408'
Set-TextValue 'D7' '6.380'
Set-TextValue 'D8' '0.8156'
Set-TextValue 'D9' '0.9265'
Set-TextValue 'B10' 'WazirX'
Set-TextValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1438'
Set-TextValue 'E10' '9WazirXWRX'
Set-TextValue 'B11' 'MandalaExchangeToken'
Set-TextValue 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.07517'
Set-TextValue 'E11' '10MandalaExchangeTokenMDX'
Set-TextValue 'B12' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C12' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D12' '0.03127'
Set-TextValue 'E12' '11LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03097'
Set-TextValue 'E13' '12BitrueCoinBTR'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09374'
Set-TextValue 'E14' '13BitMartTokenBMX'
Set-TextValue 'B15' 'MCDex'
Set-TextValue 'C15' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D15' '3.555'
Set-TextValue 'E15' '14MCDexMCB'
Set-TextValue 'B16' 'BitForexToken'
Set-TextValue 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D16' '0.001590'
Set-TextValue 'E16' '15BitForexTokenBF'
Set-TextValue 'B17' 'CoinExToken'
Set-TextValue 'C17' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D17' '0.04772'
Set-TextValue 'E17' '16CoinExTokenCET'
Set-TextValue 'B18' 'One'
Set-TextValue 'C18' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D18' '0.0005784'
Set-TextValue 'E18' '17OneONE'
Set-TextValue 'D19' '0.006394'
Set-TextValue 'D20' '0.005008'
Set-TextValue 'D21' '0.001032'
Set-TextValue 'D22' '0.0001500'
Set-TextValue 'D24' '2.185'
Set-TextValue 'D26' '0.1264'
Set-TextValue 'D28' '0.0002997'
Set-TextValue 'D40' '0.04033'
Set-TextValue 'B41' 'BKEXToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D41' '0.1073'
Set-TextValue 'E41' '40BKEXTokenBKK'
Set-TextValue 'B42' 'CEJI'
Set-TextValue 'C42' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D42' '0.002710'
Set-TextValue 'E42' '41CEJICEJI'
Set-TextValue 'B43' 'KickToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D43' '0.006770'
Set-TextValue 'E43' '42KickTokenKICK'
Set-TextValue 'D44' '0.008041'
Set-TextValue 'D45' '0.00005796'
Set-TextValue 'D47' '0.4994'
Set-TextValue 'E48' '47BOLOBOLOBestin24h'

$ws.Range("ZZ1").ClearContents()
$excel.CutCopyMode = 0
